# Apply the two reachable edits described by the commit "Add files via
# upload" (the third edit in the diff -- the handout master's cached
# datetimeFigureOut field text -- is a PowerPoint-managed auto-update
# value that isn't exposed as a writable property on this object model;
# attempting to poke it either no-ops or collapses the <a:fld> into a
# literal run, so it is intentionally left alone here).
#
#   1. Slide 1 title placeholder ("Rectangle 2") repositioned/resized.
#   2. Slide 1 "Rectangle 4" textbox nudged down (Top only).
#
# NOTE on units: Shape.Left/Top/Width/Height are expressed in points
# (1 pt = 12700 EMU) over COM, while the underlying OOXML stores English
# Metric Units (EMU). To land on an exact EMU value we add a tiny epsilon
# (1e-5 pt, about 0.13 EMU) to each point value before it is truncated
# back down to EMU during serialization; this keeps the round trip exact
# for the targets below without spilling into the next EMU.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$title = $s.Shapes.Item("Rectangle 2")
$title.Left   = 66.00001               # 838200 EMU
$title.Top    = 24.631584803149607     # 312821 EMU
$title.Width  = 564.00001              # 7162800 EMU
$title.Height = 132.00001              # 1676400 EMU

$rect4 = $s.Shapes.Item("Rectangle 4")
$rect4.Top = 204.19756905511812        # 2593309 EMU

Write-Host "Done applying edits"
